$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.279.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.631.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.11%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.629.80'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.167'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.91%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000190'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.37%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.114.19'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.20%  '
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.167.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.654.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '382.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.52%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.49%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.17%  '
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.34%  '
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +15.04%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.766.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0968'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '545.65'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.44%  '
$ws.Range('E33').Value = '  +5.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.52%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.23'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.50%  '
$ws.Range('E39').Value = '  +6.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.93%  '
$ws.Range('E41').Value = '  +4.86%  '
$ws.Range('E42').Value = '  +9.78%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.332'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '154.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('E49').Value = '  +5.47%  '
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0265'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.37%  '
